$wb = $excel.ActiveWorkbook

# --- Linear sheet updates ---
$wsLinear = $wb.Worksheets.Item("Linear")
$wsLinear.Range("B2").Value = -0.08240671282437537
$wsLinear.Range("B3").Value = 1.011442406952637
$wsLinear.Range("B4").Value = 672.6095867142686
$wsLinear.Range("B5").Value = "[1.0, 0.14886215180326035, -0.043009959636977374, 0.006215959023565756, -0.028127828499673893, -0.04054390239736608, 0.1329358357615791, 0.27599967751444077, 0.1427044024429345, -0.06412641861550496, -0.020443457737528734, -0.027662016444807588, -0.05549570559559791, 0.11930897533608907, 0.2749404621355852, 0.08374098238725364, -0.053960147948551376, -0.04874449049568199, -0.034707429109932866, -0.06734962307929077]"

# --- NonLinear sheet updates ---
$wsNonLinear = $wb.Worksheets.Item("NonLinear")
$wsNonLinear.Range("B4").Value = -0.8496618998130007
$wsNonLinear.Range("B5").Value = 0.1091323069010495
$wsNonLinear.Range("B6").Value = 681.551525894158
$wsNonLinear.Range("B7").Value = 1.215012248388281
$wsNonLinear.Range("B8").Value = -1.367955328586101
$wsNonLinear.Range("B9").Value = 664.36363089392
$wsNonLinear.Range("B10").Value = "[0.9999999999999998, 0.1461978636743438, -0.04118320354161974, 0.007998324951721817, -0.02525472064084424, -0.038623382583016765, 0.13162411367239918, 0.2712840900736729, 0.14135699678666036, -0.061327135975747606, -0.01823206566815425, -0.025864979328781255, -0.052271568652723636, 0.11814601044592526, 0.27027080033512796, 0.0808819401185827, -0.05144084731236063, -0.04727984550510143, -0.03328816437210092, -0.06390755341779279]"
